$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Distance soleil" column of Tableau2 to "Angle de rotation (deg)"
# (set via the header cell so the table definition + shared string + sheet cell
# all update consistently, matching how Excel itself handles table header renames)
$tbl = $ws.ListObjects.Item("Tableau2")
$tbl.HeaderRowRange.Cells.Item(1, 5).Value = "Angle de rotation (deg)"

# Fill in the rotation angle (axial tilt, in degrees) for every body in the table
$ws.Range("E2").Value = "'7.23"
$ws.Range("E3").Value = 0.03
$ws.Range("E4").Value = 177.36
$ws.Range("E5").Value = 23.45
$ws.Range("E6").Value = 25
$ws.Range("E7").Value = 1.304
$ws.Range("E8").Value = 27
$ws.Range("E9").Value = 98
$ws.Range("E10").Value = 28.32

# Match the final selection left by the author
$ws.Range("D14").Select() | Out-Null
